# "data preparing transformers creation"
# - Update a handful of Rho (column D) outlier values
# - Drop the now-unneeded "shape"/"type" columns (F:G) and the stray
#   numeric helper values that lived in column E
# - Reposition the scatter chart that was nudged down/right once the
#   header row lost its extra columns
# - Leave the selection on the last data point in column D

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected Rho readings
$ws.Range("D2").Value = 230
$ws.Range("D5").Value = 125
$ws.Range("D9").Value = 168

# Clear the stray helper values that lived in column E (pk 0 and pk 80 rows)
$ws.Range("E2").ClearContents()
$ws.Range("E10").ClearContents()

# Remove the "shape"/"type" header + first-row values (columns F:G)
$ws.Range("F1:G2").ClearContents()

# Move/resize the chart to its new anchored position
$co = $ws.ChartObjects(1)
$co.Left = 472.75
$co.Top = 15.75
$co.Width = 433.0625
$co.Height = 216

# Match the saved selection
$ws.Range("D16").Select()
